$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Model" -> "Type"
$ws.Range("A1").Value = "Type"

# Column G model name fix: "B36TC-Bonanza" -> "B36 Bonanza"
$ws.Range("G1").Value = "B36 Bonanza"

# Replace the Gulfstream G V (column K) with the Honda HondaJet
$ws.Range("K1").Value = "HondaJet"
$ws.Range("K2").Value = "Honda"
$ws.Range("K4").Value = 422
$ws.Range("K5").Value = 368
$ws.Range("K6").Value = 13106
$ws.Range("K7").Value = 4808
$ws.Range("K8").Value = 3267
$ws.Range("K9").Value = 1206
$ws.Range("K10").Value = 13

# Relabel the metric rows (units punctuation + "Gross Weight" -> "Max Gross")
$ws.Range("A4").Value = "Max Speed, kts"
$ws.Range("A5").Value = "Cruise Speed, kts"
$ws.Range("A6").Value = "Service Ceiling, m"
$ws.Range("A7").Value = "Max Gross, kg"
$ws.Range("A8").Value = "Empty Weight, kg"
$ws.Range("A9").Value = "Range, nm"
$ws.Range("A10").Value = "Length, m"
